# Updates the cryptocurrency price/volume table on Sheet1 with refreshed
# quotes, matching the scraped GitHub Actions commit.
# Note: several "Price" cells look like plain numbers (e.g. "0.3928"); a
# leading apostrophe is used for those so Excel stores the exact text
# (preserving trailing zeros / thousands-dot notation) instead of
# converting them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.759.15'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '1.699.59'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = "'316.29"
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").Value = "'0.3928"
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").Value = "'0.4037"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -3.00%  '
$ws.Range("D10").Value = "'54.02"
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").Value = "'0.08878"
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("D13").Value = "'7.229"
$ws.Range("E13").Value = '  -1.61%  '
$ws.Range("D14").Value = "'23.36"
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = "'8.034"
$ws.Range("E15").Value = '  +5.47%  '
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '1.697.19'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").Value = "'100.10"
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("D19").Value = "'0.07009"
$ws.Range("E19").Value = '  -0.27%  '
$ws.Range("D20").Value = "'19.60"
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = "'7.019"
$ws.Range("E21").Value = '  +1.46%  '
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = "'14.46"
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("D24").Value = '24.759.94'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").Value = "'3.256"
$ws.Range("E25").Value = '  +8.92%  '
$ws.Range("E26").Value = '  +0.79%  '
$ws.Range("D27").Value = "'22.74"
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("D28").Value = "'160.98"
$ws.Range("E28").Value = '  +0.88%  '
$ws.Range("D29").Value = "'136.31"
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("D30").Value = "'5.164"
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").Value = "'7.753"
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("D32").Value = "'0.08738"
$ws.Range("E32").Value = '  +1.99%  '
$ws.Range("E33").Value = '  -4.17%  '
$ws.Range("D34").Value = "'7.188"
$ws.Range("E34").Value = '  -4.24%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = "'1.965"
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").Value = "'0.2742"
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = "'14.37"
$ws.Range("E38").Value = '  -3.03%  '
$ws.Range("D39").Value = "'0.09180"
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("E40").Value = '  -1.76%  '
$ws.Range("D41").Value = "'1.464"
$ws.Range("E41").Value = '  -1.04%  '
$ws.Range("D42").Value = "'0.7671"
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").Value = "'15.88"
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("D44").Value = "'0.7167"
$ws.Range("E44").Value = '  -1.84%  '
$ws.Range("D45").Value = "'2.569"
$ws.Range("D46").Value = "'4.216"
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").Value = "'140.77"
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").Value = "'1.310"
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = "'90.53"
$ws.Range("E50").Value = '  +2.62%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.07977"
$ws.Range("E51").Value = '  -0.70%  '
